$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look numeric, so Excel stores them as text (matches source inlineStr cells)
$textCells = @("D2", "G2", "D3", "G3", "D4", "G4", "D5", "G5", "D6", "G6", "D7", "G7", "D8", "G8", "D9", "G9", "D10", "G10", "D11", "G11", "D12", "G12", "D13", "G13", "D14", "G14", "D15", "G15", "D16", "G16", "G17", "D18", "G18", "D19", "G19", "D20", "G20", "D21", "G21", "D22", "G22", "D23", "G23", "D24", "G24", "G25", "G26", "G27", "G28", "G29", "G30", "G31", "G32", "G33", "G34", "G35", "G36", "G37", "G38", "G39", "D40", "G40", "D41", "G41", "D42", "G42", "D43", "G43", "D44", "G44", "D45", "G45", "G46", "D47", "G47", "D48", "G48", "D49", "G49", "D50", "G50", "D51", "G51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply new cell values
$ws.Range("D2").Value = '245.55'
$ws.Range("F2").Value = '25-12-2022'
$ws.Range("G2").Value = '1'
$ws.Range("D3").Value = '21.97'
$ws.Range("F3").Value = '25-12-2022'
$ws.Range("G3").Value = '1'
$ws.Range("D4").Value = '5.404'
$ws.Range("F4").Value = '25-12-2022'
$ws.Range("G4").Value = '1'
$ws.Range("D5").Value = '0.06041'
$ws.Range("F5").Value = '25-12-2022'
$ws.Range("G5").Value = '1'
$ws.Range("D6").Value = '3.404'
$ws.Range("F6").Value = '25-12-2022'
$ws.Range("G6").Value = '1'
$ws.Range("D7").Value = '0.8133'
$ws.Range("F7").Value = '25-12-2022'
$ws.Range("G7").Value = '1'
$ws.Range("D8").Value = '0.9250'
$ws.Range("F8").Value = '25-12-2022'
$ws.Range("G8").Value = '1'
$ws.Range("D9").Value = '0.1430'
$ws.Range("F9").Value = '25-12-2022'
$ws.Range("G9").Value = '1'
$ws.Range("D10").Value = '0.07463'
$ws.Range("F10").Value = '25-12-2022'
$ws.Range("G10").Value = '1'
$ws.Range("D11").Value = '0.03335'
$ws.Range("F11").Value = '25-12-2022'
$ws.Range("G11").Value = '1'
$ws.Range("D12").Value = '0.03059'
$ws.Range("F12").Value = '25-12-2022'
$ws.Range("G12").Value = '1'
$ws.Range("D13").Value = '4.013'
$ws.Range("F13").Value = '25-12-2022'
$ws.Range("G13").Value = '1'
$ws.Range("D14").Value = '0.09375'
$ws.Range("F14").Value = '25-12-2022'
$ws.Range("G14").Value = '1'
$ws.Range("D15").Value = '0.001598'
$ws.Range("F15").Value = '25-12-2022'
$ws.Range("G15").Value = '1'
$ws.Range("D16").Value = '0.04803'
$ws.Range("F16").Value = '25-12-2022'
$ws.Range("G16").Value = '1'
$ws.Range("F17").Value = '25-12-2022'
$ws.Range("G17").Value = '1'
$ws.Range("D18").Value = '0.005529'
$ws.Range("F18").Value = '25-12-2022'
$ws.Range("G18").Value = '1'
$ws.Range("D19").Value = '0.004161'
$ws.Range("F19").Value = '25-12-2022'
$ws.Range("G19").Value = '1'
$ws.Range("D20").Value = '0.0009908'
$ws.Range("F20").Value = '25-12-2022'
$ws.Range("G20").Value = '1'
$ws.Range("D21").Value = '0.00008802'
$ws.Range("F21").Value = '25-12-2022'
$ws.Range("G21").Value = '1'
$ws.Range("D22").Value = '3.651'
$ws.Range("F22").Value = '25-12-2022'
$ws.Range("G22").Value = '1'
$ws.Range("D23").Value = '6.449'
$ws.Range("F23").Value = '25-12-2022'
$ws.Range("G23").Value = '1'
$ws.Range("D24").Value = '2.189'
$ws.Range("F24").Value = '25-12-2022'
$ws.Range("G24").Value = '1'
$ws.Range("F25").Value = '25-12-2022'
$ws.Range("G25").Value = '1'
$ws.Range("F26").Value = '25-12-2022'
$ws.Range("G26").Value = '1'
$ws.Range("F27").Value = '25-12-2022'
$ws.Range("G27").Value = '1'
$ws.Range("F28").Value = '25-12-2022'
$ws.Range("G28").Value = '1'
$ws.Range("F29").Value = '25-12-2022'
$ws.Range("G29").Value = '1'
$ws.Range("F30").Value = '25-12-2022'
$ws.Range("G30").Value = '1'
$ws.Range("F31").Value = '25-12-2022'
$ws.Range("G31").Value = '1'
$ws.Range("F32").Value = '25-12-2022'
$ws.Range("G32").Value = '1'
$ws.Range("F33").Value = '25-12-2022'
$ws.Range("G33").Value = '1'
$ws.Range("F34").Value = '25-12-2022'
$ws.Range("G34").Value = '1'
$ws.Range("F35").Value = '25-12-2022'
$ws.Range("G35").Value = '1'
$ws.Range("F36").Value = '25-12-2022'
$ws.Range("G36").Value = '1'
$ws.Range("F37").Value = '25-12-2022'
$ws.Range("G37").Value = '1'
$ws.Range("F38").Value = '25-12-2022'
$ws.Range("G38").Value = '1'
$ws.Range("F39").Value = '25-12-2022'
$ws.Range("G39").Value = '1'
$ws.Range("D40").Value = '0.03986'
$ws.Range("F40").Value = '25-12-2022'
$ws.Range("G40").Value = '1'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '0.1076'
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("F41").Value = '25-12-2022'
$ws.Range("G41").Value = '1'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = '0.002721'
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("F42").Value = '25-12-2022'
$ws.Range("G42").Value = '1'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = '0.003041'
$ws.Range("E43").Value = '42KickTokenKICK'
$ws.Range("F43").Value = '25-12-2022'
$ws.Range("G43").Value = '1'
$ws.Range("D44").Value = '0.005822'
$ws.Range("F44").Value = '25-12-2022'
$ws.Range("G44").Value = '1'
$ws.Range("D45").Value = '0.00005275'
$ws.Range("F45").Value = '25-12-2022'
$ws.Range("G45").Value = '1'
$ws.Range("F46").Value = '25-12-2022'
$ws.Range("G46").Value = '1'
$ws.Range("B47").Value = 'ACDXExchange'
$ws.Range("C47").Value = 'https://coinranking.com/coin/-y35lbZ7U+acdxexchange-acxt'
$ws.Range("D47").Value = '0.0005801'
$ws.Range("E47").Value = '46ACDXExchangeACXT'
$ws.Range("F47").Value = '25-12-2022'
$ws.Range("G47").Value = '1'
$ws.Range("B48").Value = 'CoinbaseStockToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D48").Value = '0.8868'
$ws.Range("E48").Value = '47CoinbaseStockTokenCOINBestin24h'
$ws.Range("F48").Value = '25-12-2022'
$ws.Range("G48").Value = '1'
$ws.Range("B49").Value = 'BOLO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D49").Value = '0.002527'
$ws.Range("E49").Value = '48BOLOBOLOWorstin24h'
$ws.Range("F49").Value = '25-12-2022'
$ws.Range("G49").Value = '1'
$ws.Range("B50").Value = 'CryptobidCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc'
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").Value = '49CryptobidCoinCBC'
$ws.Range("F50").Value = '25-12-2022'
$ws.Range("G50").Value = '1'
$ws.Range("B51").Value = 'SpecialPowerGold'
$ws.Range("C51").Value = 'https://coinranking.com/coin/jPTWzmsWb+specialpowergold-spg'
$ws.Range("D51").Value = '0.01010'
$ws.Range("E51").Value = '50SpecialPowerGoldSPG'
$ws.Range("F51").Value = '25-12-2022'
$ws.Range("G51").Value = '1'
